$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3373
$ws.Range("F5").Value = 219
$ws.Range("F6").Value = 4863
$ws.Range("F7").Value = 473
$ws.Range("F8").Value = 306
$ws.Range("F9").Value = 179
$ws.Range("F10").Value = 636
$ws.Range("F11").Value = 285
$ws.Range("F12").Value = 44
$ws.Range("F13").Value = 17
$ws.Range("F14").Value = 665
$ws.Range("F15").Value = 290
$ws.Range("F18").Value = 148
$ws.Range("F19").Value = 348
$ws.Range("F20").Value = 4774
$ws.Range("F21").Value = 26
$ws.Range("F22").Value = 36
$ws.Range("F24").Value = 5911
$ws.Range("F26").Value = 1201
$ws.Range("F27").Value = 247
$ws.Range("F28").Value = 679
$ws.Range("F29").Value = 4424
$ws.Range("F31").Value = 99
$ws.Range("F32").Value = 130
$ws.Range("F33").Value = 877
$ws.Range("F34").Value = 76
$ws.Range("F35").Value = 10
$ws.Range("F36").Value = 799
$ws.Range("F37").Value = 869

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 38
$ws.Range("F4").Value = 14

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1094
$ws.Range("F4").Value = 41

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1094
$ws.Range("F5").Value = 41
$ws.Range("F8").Value = 3373
$ws.Range("F9").Value = 219
$ws.Range("F10").Value = 4863
$ws.Range("F11").Value = 473
$ws.Range("F12").Value = 306
$ws.Range("F13").Value = 179
$ws.Range("F14").Value = 636
$ws.Range("F15").Value = 285
$ws.Range("F16").Value = 44
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 665
$ws.Range("F19").Value = 290
$ws.Range("F21").Value = 38
$ws.Range("F23").Value = 148
$ws.Range("F24").Value = 348
$ws.Range("F25").Value = 4774
$ws.Range("F26").Value = 26
$ws.Range("F27").Value = 36
$ws.Range("F29").Value = 5911
$ws.Range("F31").Value = 1201
$ws.Range("F32").Value = 247
$ws.Range("F33").Value = 679
$ws.Range("F34").Value = 4424
$ws.Range("F36").Value = 14
$ws.Range("F37").Value = 99
$ws.Range("F38").Value = 130
$ws.Range("F39").Value = 877
$ws.Range("F40").Value = 76
$ws.Range("F41").Value = 10
$ws.Range("F42").Value = 799
$ws.Range("F43").Value = 869
